$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the FilesTab Neo4j query text in B4: drop the "File Type" and
#     "Breed" lines from the RETURN clause (per "corrected ICDC Breed 1-14
#     scripts"), keeping the rest of the query intact. ---
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN['German Shorthaired Pointer'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# Row 4 got shorter (two fewer wrapped lines) -> smaller autofit height.
$ws.Rows(4).RowHeight = 217.5

# --- Selection / scroll position moved: the sheet view no longer has a
#     frozen/scrolled topLeftCell, and the active selection is now B4
#     (was C5, scrolled so A4 was the top-left cell). ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B4").Select()
